$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.116.64"
$ws.Range("E2").Value = "  +11.20%  "

$ws.Range("D3").Value = "3.273.61"
$ws.Range("E3").Value = "  +6.59%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'399.65"
$ws.Range("E5").Value = "  +0.89%  "

$ws.Range("D6").Value = "'110.34"
$ws.Range("E6").Value = "  +8.32%  "

$ws.Range("E7").Value = "  +4.64%  "

$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  +6.63%  "

$ws.Range("E10").Value = "  +6.60%  "

$ws.Range("D11").Value = "'0.0955"
$ws.Range("E11").Value = "  +12.38%  "

$ws.Range("E12").Value = "  +2.39%  "

$ws.Range("D13").Value = "3.777.67"
$ws.Range("E13").Value = "  +6.28%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'19.14"
$ws.Range("E14").Value = "  +4.35%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'8.10"
$ws.Range("E15").Value = "  +5.48%  "

$ws.Range("D16").Value = "3.273.72"
$ws.Range("E16").Value = "  +5.88%  "

$ws.Range("E17").Value = "  +2.75%  "

$ws.Range("E18").Value = "  +3.48%  "

$ws.Range("D19").Value = "56.894.73"
$ws.Range("E19").Value = "  +10.81%  "

$ws.Range("D20").Value = "'3.31"
$ws.Range("E20").Value = "  +4.57%  "

$ws.Range("E21").Value = "  +11.38%  "

$ws.Range("D22").Value = "'12.93"
$ws.Range("E22").Value = "  +4.97%  "

$ws.Range("D23").Value = "'304.74"
$ws.Range("E23").Value = "  +15.09%  "

$ws.Range("D24").Value = "'75.21"
$ws.Range("E24").Value = "  +7.67%  "

$ws.Range("E25").Value = "  +0.60%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'28.34"
$ws.Range("E26").Value = "  +5.44%  "

$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D27").Value = "'8.02"
$ws.Range("E27").Value = "  +1.32%  "

$ws.Range("E28").Value = "  +4.93%  "

$ws.Range("E29").Value = "  +2.39%  "

$ws.Range("E30").Value = "  +4.69%  "

$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("E32").Value = "  +4.74%  "

$ws.Range("D33").Value = "'11.05"
$ws.Range("E33").Value = "  +2.67%  "

$ws.Range("D34").Value = "'37.54"
$ws.Range("E34").Value = "  +4.43%  "

$ws.Range("D35").Value = "'0.0486"
$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("D36").Value = "'2.13"
$ws.Range("E36").Value = "  +3.21%  "

$ws.Range("D37").Value = "'51.57"
$ws.Range("E37").Value = "  +3.07%  "

$ws.Range("D38").Value = "'3.19"
$ws.Range("E38").Value = "  +26.71%  "

$ws.Range("E39").Value = "  +7.38%  "

$ws.Range("D40").Value = "'0.997"
$ws.Range("E40").Value = "  -0.27%  "

$ws.Range("E41").Value = "  +5.15%  "

$ws.Range("E42").Value = "  +5.56%  "

$ws.Range("D43").Value = "'133.94"
$ws.Range("E43").Value = "  +4.73%  "

$ws.Range("E44").Value = "  +2.08%  "

$ws.Range("E45").Value = "  +4.31%  "

$ws.Range("E46").Value = "  -3.48%  "

$ws.Range("D47").Value = "'22.16"
$ws.Range("E47").Value = "  +1.88%  "

$ws.Range("D48").Value = "2.154.36"
$ws.Range("E48").Value = "  +3.91%  "

$ws.Range("E49").Value = "  +2.39%  "

$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "'2.03"
$ws.Range("E50").Value = "  +42.45%  "

$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").Value = "'2.39"
$ws.Range("E51").Value = "  -3.52%  "
